$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Crespo record / Primera, 2022-07-27) was added
# to the daily log. It belongs chronologically right after the existing
# row 444, so insert a fresh row at 445 — this pushes the former rows
# 445-488 down to 446-489 (the workbook's dimension grows to A1:R489) and
# Excel's "insert copies formatting from the row above" behaviour carries
# the date number format (style index 2) onto the new D445 automatically.
$ws.Rows.Item(445).Insert()

$ws.Cells.Item(445, 1).Value = 4
$ws.Cells.Item(445, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(445, 3).Value = "Los Lagos"
$ws.Cells.Item(445, 4).Value = 44769
$ws.Cells.Item(445, 5).Value = 10
$ws.Cells.Item(445, 6).Value = 100112006
$ws.Cells.Item(445, 7).Value = "Repollo"
$ws.Cells.Item(445, 8).Value = "Crespo record"
$ws.Cells.Item(445, 9).Value = "Primera"
$ws.Cells.Item(445, 10).Value = 250
$ws.Cells.Item(445, 11).Value = 2000
$ws.Cells.Item(445, 12).Value = 2000
$ws.Cells.Item(445, 13).Value = 2000
$ws.Cells.Item(445, 14).Value = "$/unidad"
$ws.Cells.Item(445, 15).Value = "Región Metropolitana"
$ws.Cells.Item(445, 16).Value = 2000
$ws.Cells.Item(445, 17).Value = 1
$ws.Cells.Item(445, 18).Value = "Hortaliza"
